$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 191.66667
$ws.Range("I18").Value = 190.4762
$ws.Range("K18").Value = 190.4762
$ws.Range("M18").Value = 93.52379999999999
$ws.Range("H19").Value = 802.5454999999999
$ws.Range("I19").Value = 748.375
$ws.Range("J19").Value = 947
$ws.Range("K19").Value = 748.375
$ws.Range("L19").Value = 947
$ws.Range("M19").Value = -573.375
$ws.Range("N19").Value = -1297
$ws.Range("H141").Value = 3271.5
$ws.Range("I141").Value = 3033.3333
$ws.Range("J141").Value = 3373.5715
$ws.Range("K141").Value = 9099.999899999999
$ws.Range("L141").Value = 10120.7145
$ws.Range("M141").Value = -3919.999899999999
$ws.Range("N141").Value = -20480.7145
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 9000
$ws.Range("J17").Value = 12000
$ws.Range("K17").Value = 9000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = -8827
$ws.Range("N17").Value = -12346
$ws.Range("H18").Value = 62509.75
$ws.Range("J18").Value = 62509.75
$ws.Range("L18").Value = 62509.75
$ws.Range("N18").Value = -63153.75
$ws.Range("H23").Value = 80003.5
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H18").Value = 70011
$ws.Range("J18").Value = 70011
$ws.Range("L18").Value = 70011
$ws.Range("N18").Value = -71069
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H23").Value = 45000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 45000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 45000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -45566
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H19").Value = 1960
$ws.Range("I19").Value = 1960
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 5880
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -5706
$ws.Range("N19").ClearContents()
$ws.Range("H122").Value = 478.9375
$ws.Range("J122").Value = 821.8
$ws.Range("L122").Value = 7396.2
$ws.Range("N122").Value = -12296.2
$ws.Range("H131").Value = 16234634
$ws.Range("I131").Value = 606
$ws.Range("J131").Value = 17361998
$ws.Range("K131").Value = 1818
$ws.Range("L131").Value = 52085994
$ws.Range("M131").Value = 3222
$ws.Range("N131").Value = -52096074
$ws.Range("H138").Value = 9011755
$ws.Range("I138").Value = 1204.375
$ws.Range("J138").Value = 15876936
$ws.Range("K138").Value = 3613.125
$ws.Range("L138").Value = 47630808
$ws.Range("M138").Value = 1526.875
$ws.Range("N138").Value = -47641088
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 13165
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 13165
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 13165
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -13389
$ws.Range("H11").Value = 22457000
$ws.Range("I11").Value = 30875250
$ws.Range("J11").Value = 8333.333000000001
$ws.Range("K11").Value = 30875250
$ws.Range("L11").Value = 8333.333000000001
$ws.Range("M11").Value = -30875111
$ws.Range("N11").Value = -8611.333000000001
$ws.Range("H12").Value = 69502
$ws.Range("J12").Value = 70004
$ws.Range("L12").Value = 70004
$ws.Range("N12").Value = -70284
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("H18").Value = 70006
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 70006
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 70006
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -70592
$ws.Range("H19").Value = 10040
$ws.Range("I19").Value = 9666.666999999999
$ws.Range("J19").Value = 10600
$ws.Range("K19").Value = 9666.666999999999
$ws.Range("L19").Value = 10600
$ws.Range("M19").Value = -9378.666999999999
$ws.Range("N19").Value = -11176
$ws.Range("H23").Value = 2226.5
$ws.Range("I23").Value = 890
$ws.Range("J23").Value = 2672
$ws.Range("K23").Value = 890
$ws.Range("L23").Value = 2672
$ws.Range("M23").Value = -667
$ws.Range("N23").Value = -3118
$ws.Range("H25").Value = 47861
$ws.Range("I25").Value = 9000
$ws.Range("J25").Value = 63405.4
$ws.Range("K25").Value = 9000
$ws.Range("L25").Value = 63405.4
$ws.Range("M25").Value = -8471
$ws.Range("N25").Value = -64463.4
$ws.Range("H122").Value = 2338.606
$ws.Range("I122").Value = 2075.2666
$ws.Range("J122").Value = 4972
$ws.Range("K122").Value = 6225.7998
$ws.Range("L122").Value = 14916
$ws.Range("M122").Value = -3775.7998
$ws.Range("N122").Value = -19816
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 48105.4
$ws.Range("I23").Value = 15253
$ws.Range("J23").Value = 70007
$ws.Range("K23").Value = 15253
$ws.Range("L23").Value = 70007
$ws.Range("M23").Value = -15023
$ws.Range("N23").Value = -70467
$ws.Range("H25").Value = 24586
$ws.Range("I25").Value = 1500
$ws.Range("J25").Value = 29203.2
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 29203.2
$ws.Range("M25").Value = -1270
$ws.Range("N25").Value = -29663.2
$ws.Range("H68").Value = 2072.9524
$ws.Range("I68").Value = 1586
$ws.Range("J68").Value = 2722.2222
$ws.Range("K68").Value = 1586
$ws.Range("L68").Value = 2722.2222
$ws.Range("M68").Value = -837
$ws.Range("N68").Value = -4220.2222
$ws.Range("H71").Value = 2072.9524
$ws.Range("I71").Value = 1586
$ws.Range("J71").Value = 2722.2222
$ws.Range("K71").Value = 7930
$ws.Range("L71").Value = 13611.111
$ws.Range("M71").Value = -4186
$ws.Range("N71").Value = -21099.111
$ws.Range("H131").Value = 59999
$ws.Range("J131").Value = 59999
$ws.Range("L131").Value = 59999
$ws.Range("N131").Value = -70079
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 25845.428
$ws.Range("J6").Value = 25845.428
$ws.Range("L6").Value = 25845.428
$ws.Range("N6").Value = -26075.428
$ws.Range("H11").Value = 70005
$ws.Range("J11").Value = 70005
$ws.Range("L11").Value = 70005
$ws.Range("N11").Value = -70289
$ws.Range("H12").Value = 34574.43
$ws.Range("I12").Value = 8000
$ws.Range("K12").Value = 8000
$ws.Range("M12").Value = -7858
$ws.Range("H13").Value = 999
$ws.Range("J13").Value = 999
$ws.Range("L13").Value = 999
$ws.Range("N13").Value = -1279
$ws.Range("H17").Value = 5002500
$ws.Range("I17").Value = 5002500
$ws.Range("K17").Value = 5002500
$ws.Range("M17").Value = -5002328
$ws.Range("H18").Value = 70007
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 70007
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 70007
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -70353
$ws.Range("H19").Value = 29402.2
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 29402.2
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 29402.2
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -29750.2
$ws.Range("H23").Value = 38005.5
$ws.Range("J23").Value = 47674
$ws.Range("L23").Value = 47674
$ws.Range("N23").Value = -48132
$ws.Range("H24").Value = 70010
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 70010
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 70010
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -70470
$ws.Range("H25").Value = 60013.5
$ws.Range("J25").Value = 60013.5
$ws.Range("L25").Value = 60013.5
$ws.Range("N25").Value = -60599.5
$ws.Range("H111").Value = 53000
$ws.Range("J111").Value = 53000
$ws.Range("L111").Value = 53000
$ws.Range("N111").Value = -61180
$ws.Range("H113").Value = 396.05884
$ws.Range("I113").Value = 314.41666
$ws.Range("J113").Value = 592
$ws.Range("K113").Value = 943.2499799999999
$ws.Range("L113").Value = 1776
$ws.Range("M113").Value = 1226.75002
$ws.Range("N113").Value = -6116
